$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Assistant Professor (starting Dec 2025)", $true, $false, $false, $false, $false, $true, 1, $false, "Assistant Professor (starting Dec 2024)", 2)
